$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Pancho Jimenez
$ws.Range("A7").Value = "Pancho Jimenez"
$ws.Range("B7").Value = "12345678W"
$ws.Range("C7").Value = "arrobagmail.com"
$ws.Range("D7").Value = 246

# Row 8: Tiburcio Perez
$ws.Range("A8").Value = "Tiburcio Perez"
$ws.Range("B8").Value = "87654321W"
$ws.Range("C8").Value = "eltibu@email.com"
$ws.Range("D8").Value = 123
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:eltibu@email.com")

# Row 9: Miguel Llano
$ws.Range("A9").Value = "Miguel Llano"
$ws.Range("B9").Value = "61923982R"
$ws.Range("C9").Value = "llano@mail.com"
$ws.Range("D9").Value = 201
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:llano@mail.com")

# Row 10: Javier Ignacio Molina
$ws.Range("A10").Value = "Javier Ignacio Molina"
$ws.Range("B10").Value = "02710830G"
$ws.Range("C10").Value = "molina.com"
$ws.Range("D10").Value = 256
$ws.Range("A10").Font.Color = 0

# Row 11: Sixto Naranjo Marín
$ws.Range("A11").Value = "Sixto Naranjo Marín"
$ws.Range("B11").Value = "77631962Q"
$ws.Range("C11").Value = "sixton@email.com"
$ws.Range("D11").Value = 213
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:sixton@email.com")
$ws.Range("A11").Font.Color = 0

# Row 12: Oscar Darío Murillo
$ws.Range("A12").Value = "Oscar Darío Murillo"
$ws.Range("B12").Value = "54811130Z"
$ws.Range("C12").Value = "murillo@email.com"
$ws.Range("D12").Value = 145
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:murillo@email.com")
$ws.Range("A12").Font.Color = 0

# Row 13: Arturo Tabares
$ws.Range("A13").Value = "Arturo Tabares"
$ws.Range("B13").Value = "44788410G"
$ws.Range("C13").Value = "searturo@email.com"
$ws.Range("D13").Value = 167
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:searturo@email.com")
$ws.Range("A13").Font.Color = 0

# Row 14: Gabriel Jaime Jiménez
$ws.Range("A14").Value = "Gabriel Jaime Jiménez"
$ws.Range("B14").Value = "22965185F"
$ws.Range("C14").Value = "gbj@email.com"
$ws.Range("D14").Value = 189
$ws.Hyperlinks.Add($ws.Range("C14"), "mailto:gbj@email.com")
$ws.Range("A14").Font.Color = 0

# Row 15: Bernardo Posada Vera
$ws.Range("A15").Value = "Bernardo Posada Vera"
$ws.Range("B15").Value = "99352012Q"
$ws.Range("C15").Value = "berni@gmail.com"
$ws.Range("D15").Value = 125
$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:berni@gmail.com")
$ws.Range("A15").Font.Color = 0

# Row 16: Luis Guillermo Vélez Osorio
$ws.Range("A16").Value = "Luis Guillermo Vélez Osorio"
$ws.Range("B16").Value = "04828404Z"
$ws.Range("C16").Value = "lsgh@.com"
$ws.Range("D16").Value = 109
$ws.Hyperlinks.Add($ws.Range("C16"), "mailto:lsgh@.com")
$ws.Range("A16").Font.Color = 0

# Row 17: Horacio Augusto Moreno Correa
$ws.Range("A17").Value = "Horacio Augusto Moreno Correa"
$ws.Range("B17").Value = "56830428G"
$ws.Range("C17").Value = "a@.com.es"
$ws.Range("D17").Value = 123
$ws.Hyperlinks.Add($ws.Range("C17"), "mailto:a@.com.es")
$ws.Range("A17").Font.Color = 0

# Row 18: Alejandro Alzate Garcés
$ws.Range("A18").Value = "Alejandro Alzate Garcés"
$ws.Range("B18").Value = "66051967S"
$ws.Range("C18").Value = "alex@email.com"
$ws.Range("D18").Value = 231
$ws.Hyperlinks.Add($ws.Range("C18"), "mailto:alex@email.com")
$ws.Range("A18").Font.Color = 0

# Row 19: Gustavo Hernán Rodríguez Vallejo
$ws.Range("A19").Value = "Gustavo Hernán Rodríguez Vallejo"
$ws.Range("B19").Value = "41998336Z"
$ws.Range("C19").Value = "gsss@email.com"
$ws.Range("D19").Value = 233
$ws.Hyperlinks.Add($ws.Range("C19"), "mailto:gsss@email.com")
$ws.Range("A19").Font.Color = 0

# Row 20: Beatriz Elena Puerta
$ws.Range("A20").Value = "Beatriz Elena Puerta"
$ws.Range("B20").Value = "05164173Y"
$ws.Range("C20").Value = "puertadoor@gmail.com"
$ws.Range("D20").Value = 123
$ws.Hyperlinks.Add($ws.Range("C20"), "mailto:puertadoor@gmail.com")
$ws.Range("A20").Font.Color = 0

# Row 21: Álvaro de Jesús
$ws.Range("A21").Value = "Álvaro de Jesús"
$ws.Range("B21").Value = "10797551V"
$ws.Range("C21").Value = "jesuuuuh@mail.com"
$ws.Range("D21").Value = 187
$ws.Hyperlinks.Add($ws.Range("C21"), "mailto:jesuuuuh@mail.com")
$ws.Range("A21").Font.Color = 0

# Row 22: Héctor Darío Bermúdez
$ws.Range("A22").Value = "Héctor Darío Bermúdez"
$ws.Range("B22").Value = "92856697Q"
$ws.Range("C22").Value = "ekthor@email.com"
$ws.Range("D22").Value = 167
$ws.Hyperlinks.Add($ws.Range("C22"), "mailto:ekthor@email.com")
$ws.Range("A22").Font.Color = 0

# Row 23: Elkin Octavio Díaz
$ws.Range("A23").Value = "Elkin Octavio Díaz"
$ws.Range("B23").Value = "82430695Y"
$ws.Range("C23").Value = "octavio@gmail.com"
$ws.Range("D23").Value = 154
$ws.Hyperlinks.Add($ws.Range("C23"), "mailto:octavio@gmail.com")
$ws.Range("A23").Font.Color = 0

# Selection / view state
$ws.Range("E4").Select()
